$wb = $excel.ActiveWorkbook

# wt_gene_list: insert a new gene "ASH1" as row 6 (alphabetically between
# ASF1 and CIN5), shifting the remaining genes down by one row.
$ws1 = $wb.Worksheets.Item("wt_gene_list")
$ws1.Rows.Item(6).Insert()
$ws1.Range("A6").Value = "ASH1"

# The inserted row copied formatting from the row above; the refreshed
# gene list (A2:A17) reverts to the default/Normal style (no explicit
# font override), matching the rest of the workbook's other gene-list
# sheets.
$ws1.Range("A2:A17").Style = "Normal"

# dZAP1_gene_list was the previously active/selected sheet; clear its
# selection back to the full gene range.
$ws5 = $wb.Worksheets.Item("dZAP1_gene_list")
$ws5.Range("A2:A16").Select()

# wt_gene_list becomes the active sheet/tab, with B14 selected.
$ws1.Activate()
$ws1.Range("B14").Select()
